$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 51
$ws_ALC.Range("H51").Value = 3911.1333
$ws_ALC.Range("I51").Value = 4093.28
$ws_ALC.Range("K51").Value = 4093.28
$ws_ALC.Range("M51").Value = -3609.28

# ALC row 53
$ws_ALC.Range("H53").Value = 984.5
$ws_ALC.Range("J53").Value = 1113.8
$ws_ALC.Range("L53").Value = 1113.8
$ws_ALC.Range("N53").Value = -2387.8

# ALC row 58
$ws_ALC.Range("H58").Value = 1095
$ws_ALC.Range("I58").Value = 642.5
$ws_ALC.Range("K58").Value = 1927.5
$ws_ALC.Range("M58").Value = -1777.5

# ALC row 98
$ws_ALC.Range("H98").Value = 1143.1538
$ws_ALC.Range("J98").Value = 978
$ws_ALC.Range("L98").Value = 978
$ws_ALC.Range("N98").Value = -3974

# ALC row 113
$ws_ALC.Range("H113").Value = 3946.5715
$ws_ALC.Range("J113").Value = 4890.5
$ws_ALC.Range("L113").Value = 4890.5
$ws_ALC.Range("N113").Value = -11398.5

# ALC row 115
$ws_ALC.Range("H115").Value = 569.5
$ws_ALC.Range("I115").Value = 569.5
$ws_ALC.Range("K115").Value = 1708.5
$ws_ALC.Range("M115").Value = -141.5

# ALC row 122
$ws_ALC.Range("H122").Value = 1143.1538
$ws_ALC.Range("J122").Value = 978
$ws_ALC.Range("L122").Value = 2934
$ws_ALC.Range("N122").Value = -7834

# ALC row 132
$ws_ALC.Range("H132").Value = 22877604
$ws_ALC.Range("I132").Value = 23570776
$ws_ALC.Range("J132").Value = 2895
$ws_ALC.Range("K132").Value = 70712328
$ws_ALC.Range("L132").Value = 8685
$ws_ALC.Range("M132").Value = -70709798
$ws_ALC.Range("N132").Value = -13745

# ARM row 2
$ws_ARM.Range("H2").Value = 2009.5
$ws_ARM.Range("I2").Value = 859.8
$ws_ARM.Range("K2").Value = 859.8
$ws_ARM.Range("M2").Value = -746.8

# ARM row 4
$ws_ARM.Range("H4").Value = 366
$ws_ARM.Range("I4").Value = 207.5
$ws_ARM.Range("K4").Value = 207.5
$ws_ARM.Range("M4").Value = -91.5

# ARM row 45
$ws_ARM.Range("H45").Value = 4278.273
$ws_ARM.Range("I45").Value = 4220.2856
$ws_ARM.Range("K45").Value = 4220.2856
$ws_ARM.Range("M45").Value = -3843.2856

# ARM row 116
$ws_ARM.Range("H116").Value = 2009.5
$ws_ARM.Range("I116").Value = 859.8
$ws_ARM.Range("K116").Value = 859.8
$ws_ARM.Range("M116").Value = 1434.2

# ARM row 132
$ws_ARM.Range("H132").Value = 3724.85
$ws_ARM.Range("I132").Value = 3789.5557
$ws_ARM.Range("K132").Value = 11368.6671
$ws_ARM.Range("M132").Value = -8838.667099999999

# BSM row 3
$ws_BSM.Range("H3").Value = 2009.5
$ws_BSM.Range("I3").Value = 859.8
$ws_BSM.Range("K3").Value = 859.8
$ws_BSM.Range("M3").Value = -745.8

# BSM row 20
$ws_BSM.Range("H20").Value = 3652
$ws_BSM.Range("I20").Value = 3086
$ws_BSM.Range("K20").Value = 3086
$ws_BSM.Range("M20").Value = -2839

# BSM row 26
$ws_BSM.Range("H26").Value = 115000.75
$ws_BSM.Range("J26").Value = 120001
$ws_BSM.Range("L26").Value = 120001
$ws_BSM.Range("N26").Value = -120585

# BSM row 86
$ws_BSM.Range("H86").Value = 2773.4883
$ws_BSM.Range("I86").Value = 2643.9333
$ws_BSM.Range("J86").Value = 3072.4614
$ws_BSM.Range("K86").Value = 2643.9333
$ws_BSM.Range("L86").Value = 3072.4614
$ws_BSM.Range("M86").Value = -1520.9333
$ws_BSM.Range("N86").Value = -5318.4614

# BSM row 89
$ws_BSM.Range("H89").Value = 2773.4883
$ws_BSM.Range("I89").Value = 2643.9333
$ws_BSM.Range("J89").Value = 3072.4614
$ws_BSM.Range("K89").Value = 13219.6665
$ws_BSM.Range("L89").Value = 15362.307
$ws_BSM.Range("M89").Value = -7603.666500000001
$ws_BSM.Range("N89").Value = -26594.307

# BSM row 96
$ws_BSM.Range("H96").Value = 87721.5
$ws_BSM.Range("I96").Value = 64944
$ws_BSM.Range("J96").Value = 110499
$ws_BSM.Range("K96").Value = 64944
$ws_BSM.Range("L96").Value = 110499
$ws_BSM.Range("M96").Value = -62198
$ws_BSM.Range("N96").Value = -115991

# BSM row 134
$ws_BSM.Range("H134").Value = 4765293.5
$ws_BSM.Range("I134").Value = 5131470
$ws_BSM.Range("K134").Value = 15394410
$ws_BSM.Range("M134").Value = -15391875

# CRP row 7
$ws_CRP.Range("H7").Value = 96.64286
$ws_CRP.Range("I7").Value = 32.090908
$ws_CRP.Range("K7").Value = 32.090908
$ws_CRP.Range("M7").Value = 80.909092

# CRP row 16
$ws_CRP.Range("H16").Value = 1268.2
$ws_CRP.Range("I16").Value = 1268.2
$ws_CRP.Range("K16").Value = 1268.2
$ws_CRP.Range("M16").Value = -981.2

# CRP row 31
$ws_CRP.Range("H31").Value = 8341388
$ws_CRP.Range("J31").Value = 9405.9375
$ws_CRP.Range("L31").Value = 9405.9375
$ws_CRP.Range("N31").Value = -9995.9375

# CRP row 34
$ws_CRP.Range("H34").Value = 8341388
$ws_CRP.Range("J34").Value = 9405.9375
$ws_CRP.Range("L34").Value = 9405.9375
$ws_CRP.Range("N34").Value = -9809.9375

# CRP row 87
$ws_CRP.Range("H87").Value = 67958.8
$ws_CRP.Range("J87").Value = 67958.8
$ws_CRP.Range("L87").Value = 67958.8
$ws_CRP.Range("N87").Value = -70330.8

# CRP row 90
$ws_CRP.Range("H90").Value = 67958.8
$ws_CRP.Range("J90").Value = 67958.8
$ws_CRP.Range("L90").Value = 203876.4
$ws_CRP.Range("N90").Value = -215732.4

# CRP row 105
$ws_CRP.Range("H105").Value = 2851
$ws_CRP.Range("I105").Value = 2443.1667
$ws_CRP.Range("J105").Value = 3666.6667
$ws_CRP.Range("K105").Value = 2443.1667
$ws_CRP.Range("L105").Value = 3666.6667
$ws_CRP.Range("M105").Value = -696.1667000000002
$ws_CRP.Range("N105").Value = -7160.6667

# CRP row 113
$ws_CRP.Range("H113").Value = 1268.2
$ws_CRP.Range("I113").Value = 1268.2
$ws_CRP.Range("K113").Value = 1268.2
$ws_CRP.Range("M113").Value = 901.8

# CUL row 3
$ws_CUL.Range("H3").Value = 8142.6665
$ws_CUL.Range("I3").Value = 8142.6665
$ws_CUL.Range("K3").Value = 24427.9995
$ws_CUL.Range("M3").Value = -24315.9995

# CUL row 4
$ws_CUL.Range("H4").Value = 103515624
$ws_CUL.Range("I4").Value = 54939812
$ws_CUL.Range("J4").Value = 228424850
$ws_CUL.Range("K4").Value = 164819436
$ws_CUL.Range("L4").Value = 685274550
$ws_CUL.Range("M4").Value = -164819324
$ws_CUL.Range("N4").Value = -685274774

# CUL row 5
$ws_CUL.Range("H5").Value = 998
$ws_CUL.Range("J5").Value = 998
$ws_CUL.Range("L5").Value = 2994
$ws_CUL.Range("N5").Value = -3218

# CUL row 9
$ws_CUL.Range("H9").Value = 2154.25
$ws_CUL.Range("I9").Value = 317
$ws_CUL.Range("J9").Value = 2766.6667
$ws_CUL.Range("K9").Value = 951
$ws_CUL.Range("L9").Value = 8300.000100000001
$ws_CUL.Range("M9").Value = -727
$ws_CUL.Range("N9").Value = -8748.000100000001

# CUL row 18
$ws_CUL.Range("H18").Value = 2684.75
$ws_CUL.Range("I18").Value = 1019.75
$ws_CUL.Range("K18").Value = 3059.25
$ws_CUL.Range("M18").Value = -2890.25

# CUL row 114
$ws_CUL.Range("H114").Value = 1347.75
$ws_CUL.Range("J114").Value = 522.25
$ws_CUL.Range("L114").Value = 1566.75
$ws_CUL.Range("N114").Value = -8074.75

# CUL row 122
$ws_CUL.Range("H122").Value = 14347264
$ws_CUL.Range("I122").Value = 12820845
$ws_CUL.Range("J122").Value = 15873683
$ws_CUL.Range("K122").Value = 115387605
$ws_CUL.Range("L122").Value = 142863147
$ws_CUL.Range("M122").Value = -115385155
$ws_CUL.Range("N122").Value = -142868047

# CUL row 131
$ws_CUL.Range("H131").Value = 1604
$ws_CUL.Range("I131").Value = 1323.4
$ws_CUL.Range("K131").Value = 3970.2
$ws_CUL.Range("M131").Value = 1069.8

# CUL row 135
$ws_CUL.Range("H135").Value = 998
$ws_CUL.Range("J135").Value = 998
$ws_CUL.Range("L135").Value = 8982
$ws_CUL.Range("N135").Value = -14052

# GSM row 70
$ws_GSM.Range("H70").Value = 4550
$ws_GSM.Range("I70").Value = 4742.857
$ws_GSM.Range("K70").Value = 4742.857
$ws_GSM.Range("M70").Value = -4472.857

# GSM row 73
$ws_GSM.Range("H73").Value = 4550
$ws_GSM.Range("I73").Value = 4742.857
$ws_GSM.Range("K73").Value = 4742.857
$ws_GSM.Range("M73").Value = -3806.857

# GSM row 134
$ws_GSM.Range("H134").Value = 70790
$ws_GSM.Range("J134").Value = 70790
$ws_GSM.Range("L134").Value = 212370
$ws_GSM.Range("N134").Value = -217440

# LTW row 119
$ws_LTW.Range("H119").Value = 99891.336
$ws_LTW.Range("J119").Value = 99891.336
$ws_LTW.Range("L119").Value = 99891.336
$ws_LTW.Range("N119").Value = -109567.336

# LTW row 132
$ws_LTW.Range("H132").Value = 4865.8335
$ws_LTW.Range("I132").Value = 4795
$ws_LTW.Range("K132").Value = 14385
$ws_LTW.Range("M132").Value = -11855

# LTW row 136
$ws_LTW.Range("H136").Value = 19741.916
$ws_LTW.Range("I136").Value = 24888
$ws_LTW.Range("K136").Value = 74664
$ws_LTW.Range("M136").Value = -72114

# WVR row 132
$ws_WVR.Range("H132").Value = 2811.625
$ws_WVR.Range("I132").Value = 2505.2942
$ws_WVR.Range("J132").Value = 3555.5715
$ws_WVR.Range("K132").Value = 7515.882599999999
$ws_WVR.Range("L132").Value = 10666.7145
$ws_WVR.Range("M132").Value = -4985.882599999999
$ws_WVR.Range("N132").Value = -15726.7145
